$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Title shape ("Grupp nr: " + "11" -> single run "Grupp nr: 11")
$shpTitle = $s1.Shapes.Item(1)
$trTitle = $shpTitle.TextFrame.TextRange
$paraTitle = $trTitle.Paragraphs(1, 1)
$paraTitle.Text = "####"
$paraTitle = $trTitle.Paragraphs(1, 1)
$paraTitle.Text = "Grupp nr: 11"

# Content placeholder, 2nd paragraph ("Namn " + "och e-post till samtliga gruppmedlemmar" -> single run)
$shpContent = $s1.Shapes.Item(2)
$trContent = $shpContent.TextFrame.TextRange
$paraNamn = $trContent.Paragraphs(2, 1)
$paraNamn.Text = "####"
$paraNamn = $trContent.Paragraphs(2, 1)
$paraNamn.Text = "Namn och e-post till samtliga gruppmedlemmar"

# ---------------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shpVerktyg = $s2.Shapes.Item(2)
$trVerktyg = $shpVerktyg.TextFrame.TextRange

# 1st paragraph ("Vilka verktyg använde ni" + "?" -> single run)
$paraFraga = $trVerktyg.Paragraphs(1, 1)
$paraFraga.Text = "####"
$paraFraga = $trVerktyg.Paragraphs(1, 1)
$paraFraga.Text = "Vilka verktyg använde ni?"

# 3rd paragraph ("Ant" + " " + " " -> "Ant" + merged "  " run)
$paraAnt = $trVerktyg.Paragraphs(3, 1)
$spaces = $trVerktyg.Characters($paraAnt.Start + 3, 2)
$spaces.Text = "####"
$spaces = $trVerktyg.Characters($paraAnt.Start + 3, 4)
$spaces.Text = "  "
